# update scripts wuth new tpm
# Re-labels the "Inflammatory-Mac" cluster as "Resolving-Mac" and refreshes the
# NATMI ligand/receptor TPM-derived metrics for the Cort-Sstr3 pair table,
# adding the two new Resolving-Mac<->Resolving-Mac / Resolving-Mac<->FAPs rows
# that come with the updated cluster.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-5 (new TPM values; rename Inflammatory-Mac -> Resolving-Mac) ---
# Row 2
$ws.Range("I2").Value = 0.9158046909573684
$ws.Range("J2").Value = 0.9158046909573684
$ws.Range("O2").Value = 0.9266932386159046
$ws.Range("P2").Value = 0.9266932386159047
$ws.Range("S2").Value = 0.8486700150029214
$ws.Range("T2").Value = 0.8486700150029214

# Row 3
$ws.Range("I3").Value = 0.9158046909573684
$ws.Range("J3").Value = 0.9158046909573684
$ws.Range("O3").Value = 0.01334831496537368
$ws.Range("P3").Value = 0.01334831496537368
$ws.Range("S3").Value = 0.01222444946166565
$ws.Range("T3").Value = 0.01222444946166566

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("D4").Value = "Resolving-Mac"
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2586766666666667
$ws.Range("H4").Value = 0.77603
$ws.Range("I4").Value = 0.9158046909573684
$ws.Range("J4").Value = 0.9158046909573684
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.04007466666666667
$ws.Range("N4").Value = 0.120224
$ws.Range("O4").Value = 0.05995844641872165
$ws.Range("P4").Value = 0.05995844641872165
$ws.Range("Q4").Value = 0.01036638119111111
$ws.Range("R4").Value = 0.09329743072
$ws.Range("S4").Value = 0.05491022649278131
$ws.Range("T4").Value = 0.05491022649278132

# Row 5
$ws.Range("A5").Value = "Resolving-Mac"
$ws.Range("D5").Value = "FAPs"
$ws.Range("G5").Value = 0.02378166666666667
$ws.Range("H5").Value = 0.07134500000000001
$ws.Range("I5").Value = 0.08419530904263167
$ws.Range("J5").Value = 0.08419530904263166
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 0.6666666666666666
$ws.Range("M5").Value = 0.6193776666666667
$ws.Range("N5").Value = 1.858133
$ws.Range("O5").Value = 0.9266932386159046
$ws.Range("P5").Value = 0.9266932386159047
$ws.Range("Q5").Value = 0.01472983320944445
$ws.Range("R5").Value = 0.132568498885
$ws.Range("S5").Value = 0.07802322361298331
$ws.Range("T5").Value = 0.07802322361298329

# --- Add new rows 6-7 ---
# Row 6
$ws.Range("A6").Value = "Resolving-Mac"
$ws.Range("B6").Value = "Cort"
$ws.Range("C6").Value = "Sstr3"
$ws.Range("D6").Value = "MuSCs"
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.3333333333333333
$ws.Range("G6").Value = 0.02378166666666667
$ws.Range("H6").Value = 0.07134500000000001
$ws.Range("I6").Value = 0.08419530904263167
$ws.Range("J6").Value = 0.08419530904263166
$ws.Range("K6").Value = 1
$ws.Range("L6").Value = 0.3333333333333333
$ws.Range("M6").Value = 0.008921666666666666
$ws.Range("N6").Value = 0.026765
$ws.Range("O6").Value = 0.01334831496537368
$ws.Range("P6").Value = 0.01334831496537368
$ws.Range("Q6").Value = 0.0002121721027777778
$ws.Range("R6").Value = 0.001909548925
$ws.Range("S6").Value = 0.001123865503708022
$ws.Range("T6").Value = 0.001123865503708022

# Row 7
$ws.Range("A7").Value = "Resolving-Mac"
$ws.Range("B7").Value = "Cort"
$ws.Range("C7").Value = "Sstr3"
$ws.Range("D7").Value = "Resolving-Mac"
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.3333333333333333
$ws.Range("G7").Value = 0.02378166666666667
$ws.Range("H7").Value = 0.07134500000000001
$ws.Range("I7").Value = 0.08419530904263167
$ws.Range("J7").Value = 0.08419530904263166
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.04007466666666667
$ws.Range("N7").Value = 0.120224
$ws.Range("O7").Value = 0.05995844641872165
$ws.Range("P7").Value = 0.05995844641872165
$ws.Range("Q7").Value = 0.0009530423644444446
$ws.Range("R7").Value = 0.008577381280000001
$ws.Range("S7").Value = 0.005048219925940341
$ws.Range("T7").Value = 0.005048219925940341

Write-Output "done"